# Merge "El " + "archivo esta" + " " (with two proofErr marks in between)
# into a single run "El archivo esta " in the table cell paragraph that
# currently reads "El archivo esta vacio ".
$d = $word.ActiveDocument

$p1 = $d.Paragraphs.Item(113)
$start1 = $p1.Range.Start
$r1 = $d.Range($start1, $start1 + 16)
# First overwrite with a distinct placeholder so the engine actually
# rewrites/merges the underlying runs (a no-op same-text assignment is
# otherwise elided), then put the real text back.
$r1.Text = "El archivo esta_"
$r1b = $d.Range($start1, $start1 + 16)
$r1b.Text = "El archivo esta "

# Merge "Mismo " + "nombre" + " pero diferente " (with two proofErr marks
# in between) into a single run "Mismo nombre pero diferente " in the
# table cell paragraph that currently reads
# "Mismo nombre pero diferente extension".
$p2 = $d.Paragraphs.Item(119)
$start2 = $p2.Range.Start
$r2 = $d.Range($start2, $start2 + 28)
$r2.Text = "Mismo nombre pero diferente_"
$r2b = $d.Range($start2, $start2 + 28)
$r2b.Text = "Mismo nombre pero diferente "

# Append a new certification paragraph after the final (empty) paragraph
# of the document, right before the sectPr.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)

$newParaXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Sinespaciado"/><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Lato" w:hAnsi="Lato"/><w:color w:val="2D3B45"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t xml:space="preserve">Por medio de este p&#225;rrafo yo </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Lato" w:hAnsi="Lato"/><w:color w:val="2D3B45"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t>Roberto Cant&#250;</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Lato" w:hAnsi="Lato"/><w:color w:val="2D3B45"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t xml:space="preserve"> certifico que soy el autor intelectual del material que arriba entrego y que no realic&#233; una copia de ninguna otra persona. En aquellos casos en los que tuve que incluir material de otra persona, report&#233; su apropiada referencia.</w:t></w:r></w:p>
'@

$newRange = $d.Range($newPara.Range.Start, $newPara.Range.End)
$newRange.InsertXML($newParaXml)
